# Commit: "Fruta / hortaliza, semanal" -- weekly price update.
# A new weekly price record (week of 2022-07-07, serial 44749) is inserted
# as the new first data row for this market/category, pushing every
# existing record down by one row (old row 8 -> row 9, ..., old row 55 ->
# row 56). The sheet's used range grows from A1:R55 to A1:R56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8, shifting rows 8:55 down to 9:56.
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 with the new weekly record. The market/
# region/category/quality/unit/origin/classification columns are constant
# across this entire sheet, so reuse the same literal values as every other
# row.
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value2 = 44749
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112026
$ws.Range("G8").Value = "Haba"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 1100
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13000
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 520
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
